$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column cells to remain Text (avoid Excel auto-numeric coercion)
# by temporarily flipping NumberFormat to Text before assigning, then restoring
# the cell to the default 'Normal' style so no stray formatting is left behind.
$priceCells = @("D2", "D3", "D5", "D6", "D9", "D12", "D14", "D15", "D17", "D18", "D19", "D20", "D21", "D22", "D23", "D25", "D26", "D27", "D28", "D30", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50")
foreach ($c in $priceCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range("D2").Value = "62.892.65"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").Value = "3.450.98"
$ws.Range("E3").Value = "  -0.50%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "580.32"
$ws.Range("E5").Value = "  +0.29%  "
$ws.Range("D6").Value = "150.33"
$ws.Range("E6").Value = "  +2.18%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E8").Value = "  +1.57%  "
$ws.Range("D9").Value = "8.07"
$ws.Range("E9").Value = "  +6.22%  "
$ws.Range("E10").Value = "  -0.05%  "
$ws.Range("E11").Value = "  +4.22%  "
$ws.Range("D12").Value = "4.045.86"
$ws.Range("E12").Value = "  -0.46%  "
$ws.Range("E13").Value = "  -0.42%  "
$ws.Range("D14").Value = "28.33"
$ws.Range("E14").Value = "  -4.79%  "
$ws.Range("D15").Value = "3.445.20"
$ws.Range("E15").Value = "  -0.74%  "
$ws.Range("E16").Value = "  +1.40%  "
$ws.Range("D17").Value = "62.981.33"
$ws.Range("E17").Value = "  +0.27%  "
$ws.Range("D18").Value = "6.45"
$ws.Range("E18").Value = "  +1.92%  "
$ws.Range("D19").Value = "14.59"
$ws.Range("E19").Value = "  +1.61%  "
$ws.Range("D20").Value = "9.01"
$ws.Range("E20").Value = "  -2.10%  "
$ws.Range("D21").Value = "388.35"
$ws.Range("E21").Value = "  +0.15%  "
$ws.Range("D22").Value = "0.569"
$ws.Range("E22").Value = "  +1.64%  "
$ws.Range("D23").Value = "75.20"
$ws.Range("E23").Value = "  +0.85%  "
$ws.Range("D25").Value = "3.591.76"
$ws.Range("D26").Value = "0.0000114"
$ws.Range("E26").Value = "  -0.89%  "
$ws.Range("D27").Value = "0.186"
$ws.Range("E27").Value = "  +4.08%  "
$ws.Range("D28").Value = "7.71"
$ws.Range("E28").Value = "  +2.05%  "
$ws.Range("E29").Value = "  +0.33%  "
$ws.Range("D30").Value = "8.04"
$ws.Range("E30").Value = "  -1.48%  "
$ws.Range("E31").Value = "  +0.23%  "
$ws.Range("D32").Value = "1.00"
$ws.Range("E32").Value = "  +0.01%  "
$ws.Range("D33").Value = "1.35"
$ws.Range("E33").Value = "  -2.53%  "
$ws.Range("D34").Value = "23.33"
$ws.Range("E34").Value = "  -1.57%  "
$ws.Range("D35").Value = "5.45"
$ws.Range("E35").Value = "  +3.29%  "
$ws.Range("D36").Value = "1.65"
$ws.Range("E36").Value = "  +3.95%  "
$ws.Range("D37").Value = "31.74"
$ws.Range("E37").Value = "  +0.51%  "
$ws.Range("D38").Value = "6.97"
$ws.Range("E38").Value = "  -1.85%  "
$ws.Range("D39").Value = "169.41"
$ws.Range("E39").Value = "  +0.33%  "
$ws.Range("D40").Value = "3.487.53"
$ws.Range("D41").Value = "0.0789"
$ws.Range("E41").Value = "  +3.46%  "
$ws.Range("E42").Value = "  -1.38%  "
$ws.Range("D43").Value = "42.72"
$ws.Range("E43").Value = "  +1.10%  "
$ws.Range("D44").Value = "1.71"
$ws.Range("E44").Value = "  -0.33%  "
$ws.Range("D45").Value = "4.41"
$ws.Range("E45").Value = "  -1.69%  "
$ws.Range("D46").Value = "1.18"
$ws.Range("E46").Value = "  -1.36%  "
$ws.Range("D47").Value = "2.556.52"
$ws.Range("E47").Value = "  -1.90%  "
$ws.Range("D48").Value = "6.96"
$ws.Range("E48").Value = "  +3.11%  "
$ws.Range("D49").Value = "2.24"
$ws.Range("E49").Value = "  -0.18%  "
$ws.Range("D50").Value = "22.91"
$ws.Range("E50").Value = "  -0.29%  "
$ws.Range("E51").Value = "  +0.08%  "

foreach ($c in $priceCells) {
    $ws.Range($c).Style = "Normal"
}
